# "setup electron and switch to neon db"
#
# Content change: the Hiring Entities are renamed from their full company
# names to their short codes in the Employees sheet's hiring_entity_code
# column. "GameCove" -> "GC" and "Luxium" -> "LX". Because the Hiring
# Entities lookup sheet already used "GC"/"LX" as the Code column (and
# "GameCove Inc."/"Luxium Trading Inc." as the Name column), this rename
# collapses the now-duplicate "GC"/"LX" shared strings into one.

$wb = $excel.ActiveWorkbook

$employees = $wb.Worksheets.Item("Employees")

# hiring_entity_code column (AA) rows 2-6 were "GameCove" -> now "GC"
$employees.Range("AA2:AA6").Value = "GC"

# hiring_entity_code column (AA) rows 7-8 were "Luxium" -> now "LX"
$employees.Range("AA7:AA8").Value = "LX"

# Restore the selection/navigation state left behind after the edit.
$instructions = $wb.Worksheets.Item("Instructions")
$instructions.Activate()
$instructions.Range("A1").Select()

$departments = $wb.Worksheets.Item("Departments")
$departments.Activate()
$departments.Range("A1").Select()

$roles = $wb.Worksheets.Item("Roles")
$roles.Activate()
$roles.Range("A4").Select()

$hiringEntities = $wb.Worksheets.Item("Hiring Entities")
$hiringEntities.Activate()
$hiringEntities.Range("A1").Select()

$employees.Activate()
$employees.Range("AA8").Select()
